# expenses.xlsx edit:
#  - rename sheet "Audi Test" -> "Audi A4 B6"
#  - clear the long free-text "Description" notes (column H) for the
#    exploitation/repair rows (2-17), leaving a single blank space like the
#    already-blank rows below them
#  - remove the two trailing duplicate fuel rows (32 and 33), which also
#    drops the now-unused "2020-09-24" shared string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Audi A4 B6"

for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 8).Value = " "
}

# delete from the bottom up so row numbers of rows still to be removed
# don't shift
$ws.Rows.Item(33).Delete()
$ws.Rows.Item(32).Delete()
